$wb = $excel.ActiveWorkbook

# --- Add the last commodity row to the "wheat" sheet ---
$wheat = $wb.Worksheets.Item("wheat")
$wheat.Cells.Item(2, 1).Value = "BAT"
$wheat.Cells.Item(2, 2).Value = "Punjab"
$wheat.Cells.Item(2, 3).Value = "BMKI"
$wheat.Cells.Item(2, 4).Value = "Bihar"
$wheat.Cells.Item(2, 5).Value = "Wheat"
$wheat.Cells.Item(2, 6).Value = 1

# --- Add the last commodity row to the "rra" sheet ---
$rra = $wb.Worksheets.Item("rra")
$rra.Cells.Item(2, 1).Value = "JAT"
$rra.Cells.Item(2, 2).Value = "Jammu & Kashmir"
$rra.Cells.Item(2, 3).Value = "VSG"
$rra.Cells.Item(2, 4).Value = "Goa"
$rra.Cells.Item(2, 5).Value = "RRA"
$rra.Cells.Item(2, 6).Value = 4

# --- Add two new validation sheets at the end of the workbook: "frkcgr" and "wcgr" ---
# Both get the same header row layout used by every other sheet in this workbook.
$headers = @("From", "From State", "To", "To State", "Commodity", "Values")
$headerTemplate = $wb.Worksheets.Item("frk").Range("A1:F1")

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$frkcgr = $wb.Worksheets.Add($null, $lastSheet)
$frkcgr.Name = "frkcgr"
$frkcgr.Outline.SummaryRow = 1
$frkcgr.Outline.SummaryColumn = 1
$headerTemplate.Copy()
$frkcgr.Range("A1:F1").PasteSpecial(-4122)
for ($c = 1; $c -le 6; $c++) {
    $frkcgr.Cells.Item(1, $c).Value = $headers[$c - 1]
}

$wcgr = $wb.Worksheets.Add($null, $frkcgr)
$wcgr.Name = "wcgr"
$wcgr.Outline.SummaryRow = 1
$wcgr.Outline.SummaryColumn = 1
$headerTemplate.Copy()
$wcgr.Range("A1:F1").PasteSpecial(-4122)
for ($c = 1; $c -le 6; $c++) {
    $wcgr.Cells.Item(1, $c).Value = $headers[$c - 1]
}
